$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values that changed (simulating a "check" on the PIN details,
# introducing a couple of malformed/different PIN values for testing).
$ws.Range("A3").Value = "A012263031p"
$ws.Range("A7").Value = "A012263039Z"

# Remove the now-unneeded extra rows (8-20) that previously held repeats
# of the PIN value.
$ws.Range("A8:A20").EntireRow.Delete()

# Move the active selection to the first empty row below the data.
$ws.Range("A8").Select()
